$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = '26.848.28'
$ws.Range("E2").Value = '  -1.04%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = '1.870.83'
$ws.Range("E3").Value = '  -1.49%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("D4").Value = "'" + '1.001'
$ws.Range("E4").Value = '  -0.31%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").Value = "'" + '301.10'
$ws.Range("E5").Value = '  -1.90%  '

# Row 6: 'USDC' -> 'USDC'
$ws.Range("E6").Value = '  -0.26%  '

# Row 7: 'XRP' -> 'XRP'
$ws.Range("D7").Value = "'" + '0.5314'
$ws.Range("E7").Value = '  +1.64%  '

# Row 8: 'Cardano' -> 'Cardano'
$ws.Range("D8").Value = "'" + '0.3751'
$ws.Range("E8").Value = '  -1.41%  '

# Row 9: 'Dogecoin' -> 'Dogecoin'
$ws.Range("D9").Value = "'" + '0.07176'
$ws.Range("E9").Value = '  -1.49%  '

# Row 10: 'Solana' -> 'Solana'
$ws.Range("E10").Value = '  +0.79%  '

# Row 11: 'Polygon' -> 'Polygon'
$ws.Range("D11").Value = "'" + '0.8867'
$ws.Range("E11").Value = '  -1.75%  '

# Row 12: 'TRON' -> 'TRON'
$ws.Range("D12").Value = "'" + '0.08159'
$ws.Range("E12").Value = '  -0.15%  '

# Row 13: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D13").Value = '1.832.15'
$ws.Range("E13").Value = '  -3.09%  '

# Row 14: 'Litecoin' -> 'Litecoin'
$ws.Range("D14").Value = "'" + '93.14'
$ws.Range("E14").Value = '  -2.23%  '

# Row 15: 'Polkadot' -> 'Polkadot'
$ws.Range("D15").Value = "'" + '5.253'
$ws.Range("E15").Value = '  -1.86%  '

# Row 16: 'BinanceUSD' -> 'BinanceUSD'
$ws.Range("E16").Value = '  -0.28%  '

# Row 17: 'Avalanche' -> 'Avalanche'
$ws.Range("E17").Value = '  +0.07%  '

# Row 18: 'ShibaInu' -> 'ShibaInu'
$ws.Range("D18").Value = "'" + '0.000008533'
$ws.Range("E18").Value = '  -1.30%  '

# Row 19: 'Dai' -> 'Dai'
$ws.Range("D19").Value = "'" + '1.001'
$ws.Range("E19").Value = '  -0.14%  '

# Row 20: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D20").Value = '26.886.45'
$ws.Range("E20").Value = '  -1.04%  '

# Row 21: 'Uniswap' -> 'Uniswap'
$ws.Range("D21").Value = "'" + '4.968'
$ws.Range("E21").Value = '  -2.89%  '

# Row 22: 'Cosmos' -> 'Cosmos'
$ws.Range("E22").Value = '  -1.01%  '

# Row 23: 'Chainlink' -> 'Chainlink'
$ws.Range("D23").Value = "'" + '6.388'
$ws.Range("E23").Value = '  -1.23%  '

# Row 24: 'Monero' -> 'Monero'
$ws.Range("D24").Value = "'" + '146.94'
$ws.Range("E24").Value = '  -1.41%  '

# Row 25: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range("D25").Value = "'" + '2.256'
$ws.Range("E25").Value = '  -3.00%  '

# Row 26: 'EthereumClassic' -> 'Toncoin'
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = "'" + '1.732'
$ws.Range("E26").Value = '  -0.71%  '

# Row 27: 'Toncoin' -> 'EthereumClassic'
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = "'" + '18.02'
$ws.Range("E27").Value = '  -1.20%  '

# Row 28: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D28").Value = "'" + '114.25'
$ws.Range("E28").Value = '  -1.19%  '

# Row 29: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range("D29").Value = "'" + '4.739'
$ws.Range("E29").Value = '  -1.85%  '

# Row 30: 'Filecoin' -> 'Filecoin'
$ws.Range("D30").Value = "'" + '4.578'
$ws.Range("E30").Value = '  -6.49%  '

# Row 31: 'Stellar' -> 'Stellar'
$ws.Range("D31").Value = "'" + '0.09130'
$ws.Range("E31").Value = '  -0.87%  '

# Row 32: 'ImmutableX' -> 'ImmutableX'
$ws.Range("D32").Value = "'" + '0.7990'
$ws.Range("E32").Value = '  +0.49%  '

# Row 33: 'Hedera' -> 'Hedera'
$ws.Range("D33").Value = "'" + '0.05002'
$ws.Range("E33").Value = '  -0.76%  '

# Row 34: 'HuobiToken' -> 'HuobiToken'
$ws.Range("D34").Value = "'" + '2.984'
$ws.Range("E34").Value = '  +0.94%  '

# Row 35: 'ARBITRUM' -> 'ARBITRUM'
$ws.Range("D35").Value = "'" + '1.170'
$ws.Range("E35").Value = '  -4.10%  '

# Row 36: 'TheSandbox' -> 'TheSandbox'
$ws.Range("D36").Value = "'" + '0.6042'
$ws.Range("E36").Value = '  +5.94%  '

# Row 37: 'RenderToken' -> 'RenderToken'
$ws.Range("D37").Value = "'" + '2.595'
$ws.Range("E37").Value = '  -1.82%  '

# Row 38: 'MXToken' -> 'MXToken'
$ws.Range("D38").Value = "'" + '3.149'
$ws.Range("E38").Value = '  -6.34%  '

# Row 39: 'VeChain' -> 'VeChain'
$ws.Range("D39").Value = "'" + '0.01951'
$ws.Range("E39").Value = '  -2.21%  '

# Row 40: 'TrustWalletToken' -> 'TrustWalletToken'
$ws.Range("D40").Value = "'" + '1.071'
$ws.Range("E40").Value = '  -1.01%  '

# Row 41: 'FraxShare' -> 'FraxShare'
$ws.Range("D41").Value = "'" + '6.628'
$ws.Range("E41").Value = '  +0.63%  '

# Row 42: 'Aptos' -> 'Aptos'
$ws.Range("D42").Value = "'" + '8.857'
$ws.Range("E42").Value = '  -2.30%  '

# Row 43: 'Quant' -> 'Quant'
$ws.Range("D43").Value = "'" + '115.86'
$ws.Range("E43").Value = '  -0.44%  '

# Row 44: 'Decentraland' -> 'Decentraland'
$ws.Range("D44").Value = "'" + '0.5130'
$ws.Range("E44").Value = '  +5.10%  '

# Row 45: 'Algorand' -> 'Algorand'
$ws.Range("D45").Value = "'" + '0.1496'
$ws.Range("E45").Value = '  -1.03%  '

# Row 46: 'PaxDollar' -> 'PaxDollar'
$ws.Range("D46").Value = "'" + '1.0000'
$ws.Range("E46").Value = '  -0.32%  '

# Row 47: 'EnergySwap' -> 'EnergySwap'
$ws.Range("D47").Value = "'" + '9.901'
$ws.Range("E47").Value = '  -2.68%  '

# Row 48: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range("D48").Value = "'" + '1.624'
$ws.Range("E48").Value = '  -0.53%  '

# Row 49: 'Elrond' -> 'Elrond'
$ws.Range("D49").Value = "'" + '37.64'
$ws.Range("E49").Value = '  -1.83%  '

# Row 50: 'Cronos' -> 'Cronos'
$ws.Range("D50").Value = "'" + '0.06017'
$ws.Range("E50").Value = '  +1.01%  '

# Row 51: 'Aave' -> 'Aave'
$ws.Range("D51").Value = "'" + '62.11'
$ws.Range("E51").Value = '  -2.86%  '
